# carjacking-by-neighborhood-by-month.xlsx
# Commit: "Add data for 2021-12-19" -- the report's "through" date moves
# from December 10 to December 11, adding one more day's worth of
# carjacking counts (spread across the December column of every year,
# plus a couple of incidental updates elsewhere in the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet tab and update the matching header label -----------
$ws.Name = "Through 2021-12-11"
$ws.Range("B1").Value = "December 2021 (through December 11)"

# --- Update the affected data cells ---------------------------------------
# row 3 - North Lawndale
$ws.Range("N3").Value = 7

# row 4 - Englewood
$ws.Range("N4").Value = 2

# row 5 - West Pullman
$ws.Range("B5").Value = 2
$ws.Range("Z5").Value = 1
$ws.Range("BV5").Value = 2

# row 6 - Garfield Park
$ws.Range("AX6").Value = 4

# row 12 - Little Italy, UIC
$ws.Range("N12").Value = 1

# row 15 - Washington Heights
$ws.Range("N15").Value = 3

# row 19 - West Ridge
$ws.Range("AX19").Value = 1

# row 20 - West Loop
$ws.Range("D20").Value = 5

# row 21 - Wicker Park
$ws.Range("B21").Value = 2

# row 25 - Ashburn
$ws.Range("B25").Value = 1

# row 26 - Lake View
$ws.Range("B26").Value = 2

# row 29 - Avalon Park
$ws.Range("BJ29").Value = 1

# row 36 - Albany Park
$ws.Range("B36").Value = 1

# row 38 - Auburn Gresham
$ws.Range("B38").Value = 2

# row 39 - Brighton Park
$ws.Range("Z39").Value = 2

# row 41 - Chinatown
$ws.Range("BJ41").Value = 1

# row 51 - Loop
$ws.Range("B51").Value = 1
$ws.Range("N51").Value = 2

# row 53 - Hyde Park
$ws.Range("N53").Value = 1

# row 88 - North Center
$ws.Range("B88").Value = 1
$ws.Range("AX88").Value = 1

# row 93 - River North
$ws.Range("B93").Value = 2
